# Email for importing and exporting products
# Shifts the product rows 30-42 up by one record (row 30 takes what was in
# row 31, ..., row 41 takes what was in row 42, and row 42 wraps around to
# take what was originally in row 30), covering columns A (SKU) through
# F (Price). Column G (Status) is identical ("ON") for every row in this
# range, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 30
$lastRow = 42
$columns = @("A", "B", "C", "D", "E", "F")

# Capture the original A:F values for the affected rows before overwriting
# anything, since the shift reads from "row + 1" for every row.
$original = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = @{}
    foreach ($col in $columns) {
        $rowValues[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowValues
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($r -lt $lastRow) {
        $source = $original[$r + 1]
    } else {
        $source = $original[$firstRow]
    }
    foreach ($col in $columns) {
        $ws.Range("$col$r").Value2 = $source[$col]
    }
}
